# Trade #30 closed at 2026-02-16 21:27:43 - leadlag DOWN +0.000%
# Applies:
#  - Trade #12 (leadlag sheet row 11 / All Trades row 13) transitions OPEN -> CLOSED
#  - New Trade #30 opened (leadlag sheet row 26)
#  - Summary / leadlag / Comparison aggregate stats refreshed accordingly

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as literal text (avoids Excel auto-converting
# date-like / percent-like / numeric-like strings into numbers).
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------------------
# 1) Summary sheet - refresh aggregate rows
# ---------------------------------------------------------------------------
$sumWs = $wb.Worksheets.Item("Summary")

$sumWs.Range("C2").Value = 12
Set-TextValue $sumWs.Range("D2") "66.7%"
Set-TextValue $sumWs.Range("E2") "+2.8281%"
Set-TextValue $sumWs.Range("F2") "+0.2357%"

$sumWs.Range("C3").Value = 24
Set-TextValue $sumWs.Range("D3") "29.2%"
Set-TextValue $sumWs.Range("E3") "+2.7828%"
Set-TextValue $sumWs.Range("F3") "+0.1159%"

# ---------------------------------------------------------------------------
# 2) leadlag sheet - close Trade #12 (row 11) and append new Trade #30 (row 26)
# ---------------------------------------------------------------------------
$llWs = $wb.Worksheets.Item("leadlag")

# Trade #12 -> CLOSED
$llWs.Range("G11").Value = 69621.014042
Set-TextValue $llWs.Range("H11") "CLOSED"
$llWs.Range("I11").Value = 0.3524
$llWs.Range("J11").Value = 3.52
Set-TextValue $llWs.Range("M11") "time_exit_5min"
$llWs.Range("N11").Value = 5

# New Trade #30 (OPEN)
$llWs.Range("A26").Value = 30
Set-TextValue $llWs.Range("B26") "2026-02-16"
Set-TextValue $llWs.Range("C26") "21:27:43"
Set-TextValue $llWs.Range("D26") "leadlag"
Set-TextValue $llWs.Range("E26") "DOWN"
$llWs.Range("F26").Value = 68925.05
Set-TextValue $llWs.Range("H26") "OPEN"
$llWs.Range("I26").Value = 0
$llWs.Range("J26").Value = 0
$llWs.Range("K26").Value = 0.6959
Set-TextValue $llWs.Range("L26") "Coinbase leading with -0.070% move"
$llWs.Range("N26").Value = 0

# ---------------------------------------------------------------------------
# 3) All Trades sheet - append the now-closed Trade #12 as row 13
# ---------------------------------------------------------------------------
$atWs = $wb.Worksheets.Item("All Trades")

$atWs.Range("A13").Value = 12
Set-TextValue $atWs.Range("B13") "2026-02-16"
Set-TextValue $atWs.Range("C13") "21:22:35"
Set-TextValue $atWs.Range("D13") "leadlag"
Set-TextValue $atWs.Range("E13") "UP"
$atWs.Range("F13").Value = 69376.53
$atWs.Range("G13").Value = 69621.014042
Set-TextValue $atWs.Range("H13") "CLOSED"
$atWs.Range("I13").Value = 0.3524
$atWs.Range("J13").Value = 3.52
$atWs.Range("K13").Value = 0.6506999999999999
Set-TextValue $atWs.Range("L13") "Coinbase leading with 0.065% move"
Set-TextValue $atWs.Range("M13") "time_exit_5min"
$atWs.Range("N13").Value = 5

# ---------------------------------------------------------------------------
# 4) Comparison sheet - refresh leadlag stats row
# ---------------------------------------------------------------------------
$cmpWs = $wb.Worksheets.Item("Comparison")

$cmpWs.Range("B2").Value = 24
Set-TextValue $cmpWs.Range("C2") "29.2%"
Set-TextValue $cmpWs.Range("D2") "5.77"
Set-TextValue $cmpWs.Range("E2") "+0.4808%"
Set-TextValue $cmpWs.Range("G2") "2.47"
